$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows right after the existing "Contact" rows (currently rows 10 and 11),
# so the new rows land at 12 and 13, pushing everything below down by two rows.
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()

# Copy the formatting (style) from the existing Contact row (row 10) onto the two new rows.
$ws.Range("A10:B10").Copy()
$ws.Range("A12:B13").PasteSpecial(-4122)

# Populate the two new Contact rows with the same values as the existing Contact rows.
$ws.Range("A12").Value2 = "Contact"
$ws.Range("B12").Value2 = "No display for ContactDetail"
$ws.Range("A13").Value2 = "Contact"
$ws.Range("B13").Value2 = "No display for ContactDetail"
